# Payroll input changes: rename headers, add a TestResult column with
# Passed/Failed values, correct a few source values, drop the redundant
# "Contract ID" column, and color-code the new TestResult column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Cells.Item(1,2).Value = "Paycode"
$ws.Cells.Item(1,3).Value = "Dayfrom"
$ws.Cells.Item(1,4).Value = "Dayto"
$ws.Cells.Item(1,7).Value = "TestResult"
$ws.Cells.Item(1,8).Value = "EmpTestResult"

# --- Data corrections ---
$ws.Cells.Item(8,5).Value = "48,50"
$ws.Cells.Item(11,5).Value = "7,78"
$ws.Cells.Item(12,6).Value = "1,10"

# --- New TestResult values for each payroll-input row (2-15) ---
$ws.Cells.Item(2,7).Value = "Passed"
$ws.Cells.Item(3,7).Value = "Passed"
$ws.Cells.Item(4,7).Value = "Passed"
$ws.Cells.Item(5,7).Value = "Passed"
$ws.Cells.Item(6,7).Value = "Passed"
$ws.Cells.Item(7,7).Value = "Passed"
$ws.Cells.Item(8,7).Value = "Failed"
$ws.Cells.Item(9,7).Value = "Passed"
$ws.Cells.Item(10,7).Value = "Passed"
$ws.Cells.Item(11,7).Value = "Failed"
$ws.Cells.Item(12,7).Value = "Failed"
$ws.Cells.Item(13,7).Value = "Passed"
$ws.Cells.Item(14,7).Value = "Passed"
$ws.Cells.Item(15,7).Value = "Passed"

# --- Drop the old duplicated "Contract ID" values (column H, rows 2-15) ---
$ws.Range("H2:H15").ClearContents()

# --- Remove now-unused empty filler cells / trailing filler rows ---
$ws.Range("E23:F24").Clear()
$ws.Range("A46:H47").EntireRow.Delete()

# --- Conditional formatting: highlight TestResult column ---
$rng = $ws.Range("G1:G1048576")
$fcFailed = $rng.FormatConditions.Add(1, 3, '"Failed"')
$fcFailed.Interior.Color = 192
$fcPassed = $rng.FormatConditions.Add(1, 3, '"Passed"')
$fcPassed.Interior.Color = 5287936

Write-Output "payroll input updated"
